# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holdings detail, columns A-H)
# right before the existing "总计" (totals) sheet, and prepends a
# matching summary row to "总计".

$wb = $excel.ActiveWorkbook

# NOTE: worksheet references returned by Worksheets.Item(...) track
# *position*, not identity - once Worksheets.Add() shuffles sheet
# order/count, a previously-captured reference silently starts pointing
# at whatever sheet now sits at that old position. So every sheet handle
# below is fetched fresh, immediately before it's used, rather than
# cached across an Add() call.

# --- 1. Create the new "2022-Q1" sheet right before "总计" -----------------
$ws = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$ws.Name = "2022-Q1"

# Clone the cell formatting (bold/border/alignment header style + the
# row-index column styling) from an existing fund-holdings sheet so the
# new sheet matches the workbook's established look.
$wb.Worksheets.Item("2021-Q4").Range("A1:H10").Copy()
$wb.Worksheets.Item("2022-Q1").Range("A1:H10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'001182"
$ws.Cells.Item(2,3).Value = "易方达安心回馈混合"
$ws.Cells.Item(2,4).Value = "'90.36"
$ws.Cells.Item(2,5).Value = "'34.92"
$ws.Cells.Item(2,6).Value = "'1.62"
$ws.Cells.Item(2,7).Value = "'1.4638"
$ws.Cells.Item(2,8).Value = 6

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'003961"
$ws.Cells.Item(3,3).Value = "易方达瑞程灵活配置混合A"
$ws.Cells.Item(3,4).Value = "'34.99"
$ws.Cells.Item(3,5).Value = "'91.05"
$ws.Cells.Item(3,6).Value = "'3.52"
$ws.Cells.Item(3,7).Value = "'1.2316"
$ws.Cells.Item(3,8).Value = 6

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'002350"
$ws.Cells.Item(4,3).Value = "华安安华灵活配置混合"
$ws.Cells.Item(4,4).Value = "'42.47"
$ws.Cells.Item(4,5).Value = "'93.61"
$ws.Cells.Item(4,6).Value = "'2.42"
$ws.Cells.Item(4,7).Value = "'1.0278"
$ws.Cells.Item(4,8).Value = 6

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'003962"
$ws.Cells.Item(5,3).Value = "易方达瑞程灵活配置混合C"
$ws.Cells.Item(5,4).Value = "'9.83"
$ws.Cells.Item(5,5).Value = "'91.05"
$ws.Cells.Item(5,6).Value = "'3.52"
$ws.Cells.Item(5,7).Value = "'0.3460"
$ws.Cells.Item(5,8).Value = 6

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'163302"
$ws.Cells.Item(6,3).Value = "大摩资源优选混合(LOF)"
$ws.Cells.Item(6,4).Value = "'5.82"
$ws.Cells.Item(6,5).Value = "'81.78"
$ws.Cells.Item(6,6).Value = "'4.97"
$ws.Cells.Item(6,7).Value = "'0.2893"
$ws.Cells.Item(6,8).Value = 1

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'003839"
$ws.Cells.Item(7,3).Value = "易方达瑞通灵活配置混合A"
$ws.Cells.Item(7,4).Value = "'9.33"
$ws.Cells.Item(7,5).Value = "'32.04"
$ws.Cells.Item(7,6).Value = "'1.65"
$ws.Cells.Item(7,7).Value = "'0.1539"
$ws.Cells.Item(7,8).Value = 6

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'003882"
$ws.Cells.Item(8,3).Value = "易方达瑞弘灵活配置混合A"
$ws.Cells.Item(8,4).Value = "'6.90"
$ws.Cells.Item(8,5).Value = "'28.58"
$ws.Cells.Item(8,6).Value = "'1.36"
$ws.Cells.Item(8,7).Value = "'0.0938"
$ws.Cells.Item(8,8).Value = 8

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'003883"
$ws.Cells.Item(9,3).Value = "易方达瑞弘灵活配置混合C"
$ws.Cells.Item(9,4).Value = "'2.02"
$ws.Cells.Item(9,5).Value = "'28.58"
$ws.Cells.Item(9,6).Value = "'1.36"
$ws.Cells.Item(9,7).Value = "'0.0275"
$ws.Cells.Item(9,8).Value = 8

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'003840"
$ws.Cells.Item(10,3).Value = "易方达瑞通灵活配置混合C"
$ws.Cells.Item(10,4).Value = "'1.38"
$ws.Cells.Item(10,5).Value = "'32.04"
$ws.Cells.Item(10,6).Value = "'1.65"
$ws.Cells.Item(10,7).Value = "'0.0228"
$ws.Cells.Item(10,8).Value = 6

# --- 2. Insert a new row at the top of "总计" for the 2022-Q1 summary ------
# (Fetched fresh by name now that no further sheet Add/Delete will happen,
# so this reference stays valid for the rest of the script.)
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 9
$total.Cells.Item(2,4).Value = 4.66

# Re-number the helper index column (A) now that a row was inserted: it
# holds a simple 0-based running counter, so every existing row's index
# shifts down by one from what it used to be.
$usedRange = $total.UsedRange
$lastDataRow = $usedRange.Rows.Count
for ($r = 2; $r -le $lastDataRow; $r++) {
    $total.Cells.Item($r,1).Value = $r - 2
}
